$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -4.687379391667112
$ws.Range("C3").Value = -1.167863727494661
$ws.Range("C4").Value = -0.05450388432621493
$ws.Range("C5").Value = -0.4515867463689839
$ws.Range("C6").Value = -0.1284184568582085
